# The document's subtitle paragraph reads "(23CSE211)" - the course code's
# final digit is corrected from "1" to "2", giving "(23CSE212)".
$d = $word.ActiveDocument
$d.Content.Find.Execute("23CSE211", $true, $false, $false, $false, $false, $true, 1, $false, "23CSE212", 2)
